# Scheduled market-data refresh: updates item-price-derived columns
# (currentAveragePrice / NQ / HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ)
# for a batch of leves across the Asura data center sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62: The Mustache Suits Him
$ws.Range("H62").Value = 3340.8
$ws.Range("I62").Value = 3340.8
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 3340.8
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -2716.8
$ws.Range("N62").Value = $null

# Row 65: Forgery of Convenience (L)
$ws.Range("H65").Value = 3340.8
$ws.Range("I65").Value = 3340.8
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 16704
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -13584
$ws.Range("N65").Value = $null

# Row 129: Practical Command
$ws.Range("H129").Value = 1111.1
$ws.Range("I129").Value = 318
$ws.Range("J129").Value = 1163.9734
$ws.Range("K129").Value = 954
$ws.Range("L129").Value = 3491.9202
$ws.Range("M129").Value = 4046
$ws.Range("N129").Value = -13491.9202

# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 1999.4706
$ws.Range("I137").Value = 1355.0555
$ws.Range("J137").Value = 2724.4375
$ws.Range("K137").Value = 4065.1665
$ws.Range("L137").Value = 8173.3125
$ws.Range("M137").Value = -1515.1665
$ws.Range("N137").Value = -13273.3125

# Row 138: All-night Crafting
$ws.Range("H138").Value = 3634.0461
$ws.Range("I138").Value = 2061.7693
$ws.Range("J138").Value = 4682.231
$ws.Range("K138").Value = 6185.3079
$ws.Range("L138").Value = 14046.693
$ws.Range("M138").Value = -1045.3079
$ws.Range("N138").Value = -24326.693

$ws = $wb.Worksheets.Item("ARM")
# Row 6: Don't Hit Me One More Time
$ws.Range("H6").Value = 2000000
$ws.Range("I6").Value = 2000000
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 2000000
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -1999827

# Row 9: Headbangers' Thrall
$ws.Range("H9").Value = 37461.668
$ws.Range("I9").Value = 50000
$ws.Range("J9").Value = 34954
$ws.Range("K9").Value = 50000
$ws.Range("L9").Value = 34954
$ws.Range("M9").Value = -49830
$ws.Range("N9").Value = -35294

# Row 20: Cover Girl
$ws.Range("H20").Value = 37461.668
$ws.Range("I20").Value = 50000
$ws.Range("J20").Value = 34954
$ws.Range("K20").Value = 50000
$ws.Range("L20").Value = 34954
$ws.Range("M20").Value = -49730
$ws.Range("N20").Value = -35494

# Row 37: Get Shirty
$ws.Range("H37").Value = 18967
$ws.Range("I37").Value = 20034
$ws.Range("J37").Value = 17900
$ws.Range("K37").Value = 20034
$ws.Range("L37").Value = 17900
$ws.Range("M37").Value = -19761
$ws.Range("N37").Value = -18446

# Row 44: Very Slow Array
$ws.Range("H44").Value = 31333.334
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 31333.334
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 31333.334
$ws.Range("N44").Value = -32309.334

# Row 80: A Squire to Inspire
$ws.Range("H80").Value = 34024
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 34024
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 34024
$ws.Range("N80").Value = -36020

# Row 83: All's Fair in Highborn Assassination (L)
$ws.Range("H83").Value = 34024
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 34024
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 102072
$ws.Range("N83").Value = -112056

# Row 124: Ace of Gloves
$ws.Range("H124").Value = 44524.25
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 44524.25
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 44524.25
$ws.Range("N124").Value = -54344.25

# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 21085.166
$ws.Range("I132").Value = 24302.4
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 72907.20000000001
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -70377.20000000001
$ws.Range("N132").Value = -20057

# Row 135: Forgiveness for My Shins
$ws.Range("H135").Value = 40528
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 40528
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 40528
$ws.Range("N135").Value = -50668

$ws = $wb.Worksheets.Item("BSM")
# Row 22: Riveting Run
$ws.Range("H22").Value = 6314.375
$ws.Range("I22").Value = 8404.666999999999
$ws.Range("J22").Value = 43.5
$ws.Range("K22").Value = 8404.666999999999
$ws.Range("L22").Value = 43.5
$ws.Range("M22").Value = -8231.666999999999
$ws.Range("N22").Value = -389.5

# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 2679.3635
$ws.Range("I134").Value = 2652.3
$ws.Range("J134").Value = 2950
$ws.Range("K134").Value = 7956.900000000001
$ws.Range("L134").Value = 8850
$ws.Range("M134").Value = -5421.900000000001
$ws.Range("N134").Value = -13920

$ws = $wb.Worksheets.Item("CRP")
# Row 4: A Clogful of Camaraderie
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").Value = $null

# Row 7: Gridania's Got Talent
$ws.Range("H7").Value = 126.333336
$ws.Range("I7").Value = 90
$ws.Range("J7").Value = 199
$ws.Range("K7").Value = 90
$ws.Range("L7").Value = 199
$ws.Range("M7").Value = 23
$ws.Range("N7").Value = -425

# Row 22: Driving Up the Wall
$ws.Range("H22").Value = 5259.75
$ws.Range("I22").Value = 6120.294
$ws.Range("J22").Value = 383.33334
$ws.Range("K22").Value = 6120.294
$ws.Range("L22").Value = 383.33334
$ws.Range("M22").Value = -5770.294
$ws.Range("N22").Value = -1083.33334

# Row 124: Earring Awakening
$ws.Range("H124").Value = 55326
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 55326
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 55326
$ws.Range("N124").Value = -60236

$ws = $wb.Worksheets.Item("CUL")
# Row 68: Such a Butter Face
$ws.Range("H68").Value = 238868.58
$ws.Range("I68").Value = 303600
$ws.Range("J68").Value = 1520
$ws.Range("K68").Value = 910800
$ws.Range("L68").Value = 4560
$ws.Range("M68").Value = -909989
$ws.Range("N68").Value = -6182

# Row 71: No Margarine of Error (L)
$ws.Range("H71").Value = 238868.58
$ws.Range("I71").Value = 303600
$ws.Range("J71").Value = 1520
$ws.Range("K71").Value = 2732400
$ws.Range("L71").Value = 13680
$ws.Range("M71").Value = -2728344
$ws.Range("N71").Value = -21792

# Row 131: The Mountain Steeped
$ws.Range("H131").Value = 2404.4429
$ws.Range("I131").Value = 487.33334
$ws.Range("J131").Value = 3068.0576
$ws.Range("K131").Value = 1462.00002
$ws.Range("L131").Value = 9204.1728
$ws.Range("M131").Value = 3577.99998
$ws.Range("N131").Value = -19284.1728

# Row 133: Friends Are Food
$ws.Range("H133").Value = 5206.25
$ws.Range("I133").Value = 2797.5
$ws.Range("J133").Value = 5688
$ws.Range("K133").Value = 8392.5
$ws.Range("L133").Value = 17064
$ws.Range("M133").Value = -3332.5
$ws.Range("N133").Value = -27184

$ws = $wb.Worksheets.Item("GSM")
# Row 2: Copper and Robbers
$ws.Range("H2").Value = 41.75
$ws.Range("I2").Value = 38.166668
$ws.Range("J2").Value = 52.5
$ws.Range("K2").Value = 38.166668
$ws.Range("L2").Value = 52.5
$ws.Range("M2").Value = 74.833332
$ws.Range("N2").Value = -278.5

# Row 123: Workplace Workout
$ws.Range("H123").Value = 20999.564
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 20999.564
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 20999.564
$ws.Range("N123").Value = -25899.564

$ws = $wb.Worksheets.Item("LTW")
# Row 2: Red in the Head
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").Value = $null

$ws = $wb.Worksheets.Item("WVR")
# Row 52: Party Animals
$ws.Range("H52").Value = 3060
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 3060
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 3060
$ws.Range("M52").Value = $null
$ws.Range("N52").Value = -3512

# Row 113: A Tender Table
$ws.Range("H113").Value = 1229.625
$ws.Range("I113").Value = 950.53845
$ws.Range("J113").Value = 1559.4546
$ws.Range("K113").Value = 2851.61535
$ws.Range("L113").Value = 4678.3638
$ws.Range("M113").Value = -681.61535
$ws.Range("N113").Value = -9018.363799999999

# Row 122: Heavy Armoire
$ws.Range("H122").Value = 1387.4706
$ws.Range("I122").Value = 1160.1538
$ws.Range("J122").Value = 2126.25
$ws.Range("K122").Value = 3480.4614
$ws.Range("L122").Value = 6378.75
$ws.Range("M122").Value = -1030.4614
$ws.Range("N122").Value = -11278.75

# Row 136: Weaving the Envelope
$ws.Range("H136").Value = 2516.6843
$ws.Range("I136").Value = 2718.1667
$ws.Range("J136").Value = 2171.2856
$ws.Range("K136").Value = 8154.500100000001
$ws.Range("L136").Value = 6513.8568
$ws.Range("M136").Value = -5604.500100000001
$ws.Range("N136").Value = -11613.8568

Write-Output "Asura_Profits sheets refreshed."
